$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, pushing the old rows 8-10 down to 9-11.
$ws.Rows("8:8").Insert()

# Populate the new row 8 with the new data entry.
$ws.Cells.Item(8, 1).Value  = 10
$ws.Cells.Item(8, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value  = "La Araucanía"
$ws.Cells.Item(8, 4).Value  = 45033
$ws.Cells.Item(8, 5).Value  = 9
$ws.Cells.Item(8, 6).Value  = 100112041
$ws.Cells.Item(8, 7).Value  = "Fruto del paraíso"
$ws.Cells.Item(8, 8).Value  = "Sin especificar"
$ws.Cells.Item(8, 9).Value  = "Primera"
$ws.Cells.Item(8, 10).Value = 80
$ws.Cells.Item(8, 11).Value = 24000
$ws.Cells.Item(8, 12).Value = 24000
$ws.Cells.Item(8, 13).Value = 24000
$ws.Cells.Item(8, 14).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 1333
$ws.Cells.Item(8, 17).Value = 18
$ws.Cells.Item(8, 18).Value = "Hortaliza"

$wb.Save()
